$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 to hold the new "M1" target-cluster record.
# This pushes the existing M2 / Neutro / sCs rows down by one.
$ws.Rows.Item(4).Insert()

# Row 2: sCs -> Il1rapl1 -> Ptprf -> ECs
$ws.Cells.Item(2,1).Value = "sCs"
$ws.Cells.Item(2,2).Value = "Il1rapl1"
$ws.Cells.Item(2,3).Value = "Ptprf"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.03637433333333333
$ws.Cells.Item(2,8).Value = 0.109123
$ws.Cells.Item(2,9).Value = 1
$ws.Cells.Item(2,10).Value = 1
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.2138683333333334
$ws.Cells.Item(2,14).Value = 0.6416050000000001
$ws.Cells.Item(2,15).Value = 0.02663324640928231
$ws.Cells.Item(2,16).Value = 0.02663324640928232
$ws.Cells.Item(2,17).Value = 0.007779318046111112
$ws.Cells.Item(2,18).Value = 0.070013862415
$ws.Cells.Item(2,19).Value = 0.02663324640928231
$ws.Cells.Item(2,20).Value = 0.02663324640928232

# Row 3: sCs -> Il1rapl1 -> Ptprf -> FAPs
$ws.Cells.Item(3,1).Value = "sCs"
$ws.Cells.Item(3,2).Value = "Il1rapl1"
$ws.Cells.Item(3,3).Value = "Ptprf"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.03637433333333333
$ws.Cells.Item(3,8).Value = 0.109123
$ws.Cells.Item(3,9).Value = 1
$ws.Cells.Item(3,10).Value = 1
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.876179
$ws.Cells.Item(3,14).Value = 11.628537
$ws.Cells.Item(3,15).Value = 0.4827046100021922
$ws.Cells.Item(3,16).Value = 0.4827046100021922
$ws.Cells.Item(3,17).Value = 0.1409934270056667
$ws.Cells.Item(3,18).Value = 1.268940843051
$ws.Cells.Item(3,19).Value = 0.4827046100021922
$ws.Cells.Item(3,20).Value = 0.4827046100021922

# Row 4 (new): sCs -> Il1rapl1 -> Ptprf -> M1
$ws.Cells.Item(4,1).Value = "sCs"
$ws.Cells.Item(4,2).Value = "Il1rapl1"
$ws.Cells.Item(4,3).Value = "Ptprf"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.03637433333333333
$ws.Cells.Item(4,8).Value = 0.109123
$ws.Cells.Item(4,9).Value = 1
$ws.Cells.Item(4,10).Value = 1
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.02104366666666667
$ws.Cells.Item(4,14).Value = 0.06313100000000001
$ws.Cells.Item(4,15).Value = 0.00262058973833496
$ws.Cells.Item(4,16).Value = 0.00262058973833496
$ws.Cells.Item(4,17).Value = 0.000765449345888889
$ws.Cells.Item(4,18).Value = 0.006889044113
$ws.Cells.Item(4,19).Value = 0.00262058973833496
$ws.Cells.Item(4,20).Value = 0.00262058973833496

# Row 5: sCs -> Il1rapl1 -> Ptprf -> M2
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Il1rapl1"
$ws.Cells.Item(5,3).Value = "Ptprf"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.03637433333333333
$ws.Cells.Item(5,8).Value = 0.109123
$ws.Cells.Item(5,9).Value = 1
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.237489
$ws.Cells.Item(5,14).Value = 0.7124669999999999
$ws.Cells.Item(5,15).Value = 0.02957475264295343
$ws.Cells.Item(5,16).Value = 0.02957475264295343
$ws.Cells.Item(5,17).Value = 0.008638504048999998
$ws.Cells.Item(5,18).Value = 0.07774653644099998
$ws.Cells.Item(5,19).Value = 0.02957475264295343
$ws.Cells.Item(5,20).Value = 0.02957475264295343

# Row 6: sCs -> Il1rapl1 -> Ptprf -> Neutro
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Il1rapl1"
$ws.Cells.Item(6,3).Value = "Ptprf"
$ws.Cells.Item(6,4).Value = "Neutro"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.03637433333333333
$ws.Cells.Item(6,8).Value = 0.109123
$ws.Cells.Item(6,9).Value = 1
$ws.Cells.Item(6,10).Value = 1
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.05256933333333333
$ws.Cells.Item(6,14).Value = 0.157708
$ws.Cells.Item(6,15).Value = 0.006546513859329486
$ws.Cells.Item(6,16).Value = 0.006546513859329486
$ws.Cells.Item(6,17).Value = 0.001912174453777778
$ws.Cells.Item(6,18).Value = 0.017209570084
$ws.Cells.Item(6,19).Value = 0.006546513859329486
$ws.Cells.Item(6,20).Value = 0.006546513859329486

# Row 7: sCs -> Il1rapl1 -> Ptprf -> sCs
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Il1rapl1"
$ws.Cells.Item(7,3).Value = "Ptprf"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.03637433333333333
$ws.Cells.Item(7,8).Value = 0.109123
$ws.Cells.Item(7,9).Value = 1
$ws.Cells.Item(7,10).Value = 1
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.628977
$ws.Cells.Item(7,14).Value = 10.886931
$ws.Cells.Item(7,15).Value = 0.4519202873479077
$ws.Cells.Item(7,16).Value = 0.4519202873479077
$ws.Cells.Item(7,17).Value = 0.132001619057
$ws.Cells.Item(7,18).Value = 1.188014571513
$ws.Cells.Item(7,19).Value = 0.4519202873479077
$ws.Cells.Item(7,20).Value = 0.4519202873479077

$wb.Save()
